# Supplementary Table 1 Cultivars - add Panel/Collection info for each genotype
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Panel (column B) values for rows 1 (header) .. 37, in genotype row order
$panel = @{
    1  = "Panel"
    2  = "UNIBO"     # ACADUR
    3  = "CREA"      # ANTALIS
    4  = "GPDUR"     # ANVERGUR
    5  = "UNIBO"     # ARDENTE
    6  = "CREA"      # ASTERIX
    7  = "GPDUR"     # AVENTUR
    8  = "UNIBO"     # AZEGHAR-2_DP128
    9  = "CREA"      # BALSAMO
    10 = "UNIBO"     # BOLIDO-DP034
    11 = "CREA"      # CASANOVA
    12 = "UNIBO"     # CHAM-1_DP136
    13 = "UNIBO"     # COLOSSEO_DP087
    14 = "EPO"       # EL4X_120
    15 = "EPO"       # EL4X_194
    16 = "EPO"       # EL4X_316
    17 = "EPO"       # EL4X_35
    18 = "EPO"       # EL4X_428
    19 = "EPO"       # EL4X_464
    20 = "CREA"      # FURIO_CAMILLO
    21 = "CREA"      # GIUSTO
    22 = "UNIBO"     # KOFA
    23 = "CREA"      # L2574
    24 = "GPDUR"     # LAHAN
    25 = "GPDUR"     # LGBORIS
    26 = "UNIBO"     # LLOYD
    27 = "UNIBO"     # MIKI-1_DP161
    28 = "GPDUR"     # MONASTIR
    29 = "GPDUR"     # MURANO
    30 = "GPDUR"     # NEMESIS
    31 = "GPDUR"     # NOBILIS
    32 = "GPDUR"     # ORJAUNE
    33 = "GPDUR"     # PLUSSUR
    34 = "GPDUR"     # QUALIDOU
    35 = "CREA"      # RAMIREZ
    36 = "CREA"      # Selcuklu-97
    37 = "CREA"      # SVEVO
}

# Collection (column C) values for rows 1 (header) .. 37, in genotype row order
$collection = @{
    1  = "Collection"
    2  = "Bologna University"    # ACADUR
    3  = "CREA"                  # ANTALIS
    4  = "Bologna University"    # ANVERGUR
    5  = "Bologna University"    # ARDENTE
    6  = "CREA"                  # ASTERIX
    7  = "Arvalis"               # AVENTUR
    8  = "Bologna University"    # AZEGHAR-2_DP128
    9  = "CREA"                  # BALSAMO
    10 = "Bologna University"    # BOLIDO-DP034
    11 = "CREA"                  # CASANOVA
    12 = "Bologna University"    # CHAM-1_DP136
    13 = "Bologna University"    # COLOSSEO_DP087
    14 = "INRAE Montpellier"     # EL4X_120
    15 = "INRAE Montpellier"     # EL4X_194
    16 = "INRAE Montpellier"     # EL4X_316
    17 = "INRAE Montpellier"     # EL4X_35
    18 = "INRAE Montpellier"     # EL4X_428
    19 = "INRAE Montpellier"     # EL4X_464
    20 = "CREA"                  # FURIO_CAMILLO
    21 = "CREA"                  # GIUSTO
    22 = "Bologna University"    # KOFA
    23 = "CREA"                  # L2574
    24 = "Arvalis"               # LAHAN
    25 = "Arvalis"               # LGBORIS
    26 = "Bologna University"    # LLOYD
    27 = "Bologna University"    # MIKI-1_DP161
    28 = "Arvalis"               # MONASTIR
    29 = "Arvalis"               # MURANO
    30 = "Arvalis"               # NEMESIS
    31 = "Arvalis"               # NOBILIS
    32 = "Arvalis"               # ORJAUNE
    33 = "Arvalis"               # PLUSSUR
    34 = "Arvalis"               # QUALIDOU
    35 = "CREA"                  # RAMIREZ
    36 = "CREA"                  # Selcuklu-97
    37 = "CREA"                  # SVEVO
}

# Fill column B top to bottom (header + all data rows)
for ($r = 1; $r -le 37; $r++) {
    $ws.Cells.Item($r, 2).Value = $panel[$r]
}

# Fill column C: row 2 first, then the header (row 1), then the remaining rows
$ws.Cells.Item(2, 3).Value = $collection[2]
$ws.Cells.Item(1, 3).Value = $collection[1]
for ($r = 3; $r -le 37; $r++) {
    $ws.Cells.Item($r, 3).Value = $collection[$r]
}

# Worksheet view: zoom + selection on column D (whole column)
$ws.Activate()
$excel.ActiveWindow.Zoom = 190
$ws.Range("D1:D1048576").Select() | Out-Null

# Column widths (closest achievable values given the runtime's width quantization)
$ws.Columns.Item(1).ColumnWidth = 22.25
$ws.Columns.Item(2).ColumnWidth = 11.92
$ws.Columns.Item(3).ColumnWidth = 17.75

# Re-add the filter-database defined name (leftover from a since-removed AutoFilter),
# scoped to the sheet and hidden, matching the saved workbook metadata.
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Feuil1!`$A`$1:`$C`$37")
$fdb.Visible = $false
